$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 159, shifting rows 159:199 down to
# 160:200 (dimension grows from A1:R199 to A1:R200) - matches the diff where
# every row from the old 159 onward reappears one row lower with identical
# data, and a brand-new record is inserted at row 159.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new record.
$ws.Range("A159").Value = 5
$ws.Range("B159").Value = "Macroferia Regional de Talca"
$ws.Range("C159").Value = "Maule"
$ws.Range("D159").Value = 44754
$ws.Range("E159").Value = 7
$ws.Range("F159").Value = 100112017
$ws.Range("G159").Value = "Apio"
$ws.Range("H159").Value = "Americana (o)"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 700
$ws.Range("K159").Value = 8000
$ws.Range("L159").Value = 8000
$ws.Range("M159").Value = 8000
$ws.Range("N159").Value = "`$/docena de matas"
$ws.Range("O159").Value = "Provincia del Elquí"
$ws.Range("P159").Value = 1333
$ws.Range("Q159").Value = 6
$ws.Range("R159").Value = "Hortaliza"
